{"js": "// Add the \"Acuerdo para llamar a los barrios.\" note as a new block at the\n// end of the document, surrounded by blank paragraphs (matching the\n// formatting - bold=false/bCs=false - of the paragraph it follows).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document (ends with \"...fecha_nto_\").\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert, in order, after the last paragraph:\n//   1. an empty paragraph\n//   2. \"Acuerdo para llamar a los barrios. \"\n//   3. an empty paragraph\n//   4. an empty paragraph\n// Each insertParagraph inherits the formatting (non-bold) of the paragraph\n// it is inserted after, matching the target markup.\nconst blank1 = lastParagraph.insertParagraph(\"\", \"After\");\nconst note = blank1.insertParagraph(\"Acuerdo para llamar a los barrios. \", \"After\");\nconst blank2 = note.insertParagraph(\"\", \"After\");\nconst blank3 = blank2.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add the \"Acuerdo para llamar a los barrios.\" note as a new block at the\n# end of the document, surrounded by blank paragraphs (matching the\n# formatting - bold=false/bCs=false - of the paragraph it follows).\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document (ends with \"...fecha_nto_\").\n$lastParagraph = $d.Paragraphs.Last\n\n# Insert, in order, after the last paragraph:\n#   1. an empty paragraph\n#   2. an empty paragraph that will receive the note text\n#   3. an empty paragraph\n#   4. an empty paragraph\n# Each InsertParagraphAfter() inherits the formatting (non-bold) of the\n# paragraph it is inserted after, matching the target markup.\n$lastParagraph.Range.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$noteParagraph = $d.Paragraphs.Last\n$noteParagraph.Range.Text = \"Acuerdo para llamar a los barrios. \"\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n"}
